$d = $word.ActiveDocument

$replacements = @(
    @{old="553÷7="; new="676÷9="},
    @{old="888÷9="; new="892÷6="},
    @{old="893÷6="; new="167÷5="},
    @{old="626÷6="; new="406÷3="},
    @{old="445÷3="; new="106÷9="},
    @{old="226÷7="; new="711÷6="},
    @{old="748÷2="; new="100÷4="},
    @{old="716÷9="; new="410÷3="},
    @{old="160÷8="; new="586÷8="},
    @{old="978÷8="; new="721÷6="},
    @{old="753÷8="; new="677÷3="},
    @{old="354÷2="; new="509÷4="},
    @{old="892÷3="; new="823÷2="},
    @{old="588÷8="; new="237÷3="},
    @{old="791÷5="; new="157÷9="},
    @{old="699÷4="; new="611÷4="},
    @{old="778÷9="; new="812÷8="},
    @{old="804÷4="; new="650÷5="},
    @{old="172÷8="; new="976÷9="},
    @{old="767÷6="; new="577÷9="},
    @{old="453÷7="; new="141÷5="},
    @{old="949÷8="; new="485÷3="},
    @{old="850÷6="; new="188÷9="},
    @{old="643÷2="; new="429÷8="},
    @{old="314÷8="; new="479÷5="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
